$d = $word.ActiveDocument

# --- 1. Remove the "Meta description: ..." paragraph that follows the
#        "Play Dragon Spin free and enjoy the legendary wins" Heading1
#        paragraph. ---
$metaParaIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Meta description")) {
        $metaParaIndex = $i
        break
    }
}
if ($metaParaIndex -ne $null) {
    $metaPara = $d.Paragraphs($metaParaIndex)
    $delRng = $d.Range($metaPara.Range.Start, $metaPara.Range.End)
    $delRng.Delete()
}

# --- 2. Replace the final "Prompt: ..." paragraph with two paragraphs:
#        a new bold heading-like paragraph re-stating the page title, and
#        the former meta-description copy (minus the "Meta description: "
#        label) in italics. ---
$promptParaIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Prompt:")) {
        $promptParaIndex = $i
        break
    }
}

if ($promptParaIndex -ne $null) {
    $promptPara = $d.Paragraphs($promptParaIndex)
    $insRng = $promptPara.Range
    $insRng.Collapse(1)

    $w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
    $newXml = "<w:p $w><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Dragon Spin free and enjoy the legendary wins</w:t></w:r></w:p>" +
              "<w:p $w><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Dragon Spin online slot game and play for free. Enjoy oriental graphics and bonus rounds!</w:t></w:r></w:p>"

    $insRng.InsertXML($newXml)
}
